$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns (mirror the bold/centered/bordered style used by the
# existing header row A1:E1)
$ws.Range("F1").Value = "Linked_Posted"
$ws.Range("G1").Value = "Resume_received"
$ws.Range("H1").Value = "Resume_downloaded"

$headerRange = $ws.Range("F1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# New job posting row (Job_Id = 4)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Demo Demo"
$ws.Range("C5").Value = "DemoDemo"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
